$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column B ("R40" rule-name cell) is updated to the text value "1".
# We go through a formula -> copy -> paste-special(values) round trip rather
# than a plain `.Value = "1"` assignment so that Excel stores the result as
# genuine text (a shared string) instead of auto-coercing the numeric-looking
# string to a number, while keeping the cell's existing style/format intact
# (no quote-prefix, no new number format gets attached to the cell).
$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)
